$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values as part of repull/recalculation of data.
$updates = @{
    5  = -6
    6  = -9
    7  = -5
    10 = 9
    15 = -6
    17 = -5
    18 = -9
    19 = -9
    20 = -7
    23 = -1
    28 = -2
    32 = 1
    33 = 7
    35 = -7
    38 = 1
    39 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
